# "wrapping up test file audit"
#
# 1) Remove the stray leftover "Sheet | 3 | 4" audit row from the
#    optimization_parameters sheet (row 16). Everything below it
#    (the simulation_timepoints row, etc.) shifts up by one row, and
#    the now-unused "Sheet" shared string disappears on save.
# 2) Leave the workbook with threshold_b as the active/selected sheet
#    instead of dfhl1_log2_expression.

$wb = $excel.ActiveWorkbook

$wsOpt = $wb.Worksheets.Item("optimization_parameters")
$wsOpt.Rows.Item(16).Delete()

$wsThreshold = $wb.Worksheets.Item("threshold_b")
$wsThreshold.Activate()
